{"js": "// Replace the 25 \"NNN\u00f7N=\" division problems in the practice table with\n// their new values, per the commit's regenerated worksheet numbers.\n// Each old value is unique in the document, so an exact, case-sensitive\n// search-and-replace on each pair is unambiguous and preserves the\n// existing run formatting (font, size, etc.) of the matched text.\nconst replacements = [\n  [\"103\u00f77=\", \"938\u00f74=\"],\n  [\"401\u00f79=\", \"948\u00f74=\"],\n  [\"454\u00f76=\", \"318\u00f74=\"],\n  [\"219\u00f74=\", \"910\u00f79=\"],\n  [\"539\u00f78=\", \"298\u00f73=\"],\n  [\"588\u00f79=\", \"861\u00f74=\"],\n  [\"991\u00f74=\", \"853\u00f74=\"],\n  [\"892\u00f77=\", \"500\u00f79=\"],\n  [\"702\u00f75=\", \"817\u00f75=\"],\n  [\"712\u00f78=\", \"927\u00f72=\"],\n  [\"260\u00f73=\", \"147\u00f78=\"],\n  [\"377\u00f72=\", \"375\u00f73=\"],\n  [\"299\u00f75=\", \"260\u00f72=\"],\n  [\"867\u00f74=\", \"608\u00f74=\"],\n  [\"790\u00f72=\", \"227\u00f77=\"],\n  [\"296\u00f72=\", \"545\u00f76=\"],\n  [\"658\u00f78=\", \"644\u00f78=\"],\n  [\"963\u00f76=\", \"773\u00f76=\"],\n  [\"770\u00f77=\", \"602\u00f72=\"],\n  [\"869\u00f74=\", \"784\u00f74=\"],\n  [\"684\u00f74=\", \"907\u00f75=\"],\n  [\"692\u00f74=\", \"734\u00f74=\"],\n  [\"664\u00f76=\", \"491\u00f73=\"],\n  [\"713\u00f74=\", \"561\u00f72=\"],\n  [\"934\u00f76=\", \"441\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NNN\u00f7N=\" division problems in the practice table with\n# their new values, per the commit's regenerated worksheet numbers.\n# Each old value is unique in the document, so an exact Find/Replace on\n# each pair is unambiguous and preserves the existing run formatting\n# (font, size, etc.) of the matched text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"103\u00f77=\"; New=\"938\u00f74=\"},\n    @{Old=\"401\u00f79=\"; New=\"948\u00f74=\"},\n    @{Old=\"454\u00f76=\"; New=\"318\u00f74=\"},\n    @{Old=\"219\u00f74=\"; New=\"910\u00f79=\"},\n    @{Old=\"539\u00f78=\"; New=\"298\u00f73=\"},\n    @{Old=\"588\u00f79=\"; New=\"861\u00f74=\"},\n    @{Old=\"991\u00f74=\"; New=\"853\u00f74=\"},\n    @{Old=\"892\u00f77=\"; New=\"500\u00f79=\"},\n    @{Old=\"702\u00f75=\"; New=\"817\u00f75=\"},\n    @{Old=\"712\u00f78=\"; New=\"927\u00f72=\"},\n    @{Old=\"260\u00f73=\"; New=\"147\u00f78=\"},\n    @{Old=\"377\u00f72=\"; New=\"375\u00f73=\"},\n    @{Old=\"299\u00f75=\"; New=\"260\u00f72=\"},\n    @{Old=\"867\u00f74=\"; New=\"608\u00f74=\"},\n    @{Old=\"790\u00f72=\"; New=\"227\u00f77=\"},\n    @{Old=\"296\u00f72=\"; New=\"545\u00f76=\"},\n    @{Old=\"658\u00f78=\"; New=\"644\u00f78=\"},\n    @{Old=\"963\u00f76=\"; New=\"773\u00f76=\"},\n    @{Old=\"770\u00f77=\"; New=\"602\u00f72=\"},\n    @{Old=\"869\u00f74=\"; New=\"784\u00f74=\"},\n    @{Old=\"684\u00f74=\"; New=\"907\u00f75=\"},\n    @{Old=\"692\u00f74=\"; New=\"734\u00f74=\"},\n    @{Old=\"664\u00f76=\"; New=\"491\u00f73=\"},\n    @{Old=\"713\u00f74=\"; New=\"561\u00f72=\"},\n    @{Old=\"934\u00f76=\"; New=\"441\u00f73=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceOne = 1 (we only want the single unique match)\n    $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 1)\n}\n"}
